$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 474.5
$ws.Range("I43").Value = 450
$ws.Range("K43").Value = 450
$ws.Range("M43").Value = -381
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H94").Value = 2221
$ws.Range("I94").Value = 2221
$ws.Range("K94").Value = 2221
$ws.Range("M94").Value = -1770
$ws.Range("H121").Value = 10441.167
$ws.Range("J121").Value = 12379.4
$ws.Range("L121").Value = 37138.2
$ws.Range("N121").Value = -40632.2
$ws.Range("H137").Value = 1828.1428
$ws.Range("I137").Value = 1887.5
$ws.Range("J137").Value = 1727.6923
$ws.Range("K137").Value = 5662.5
$ws.Range("L137").Value = 5183.0769
$ws.Range("M137").Value = -3112.5
$ws.Range("N137").Value = -10283.0769
$ws.Range("H141").Value = 1135.8793
$ws.Range("I141").Value = 691.54
$ws.Range("J141").Value = 3913
$ws.Range("K141").Value = 2074.62
$ws.Range("L141").Value = 11739
$ws.Range("M141").Value = 3105.38
$ws.Range("N141").Value = -22099

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4168
$ws.Range("I32").Value = 2750.766
$ws.Range("K32").Value = 2750.766
$ws.Range("M32").Value = -2463.766
$ws.Range("H132").Value = 10989.538
$ws.Range("I132").Value = 1370.1904
$ws.Range("K132").Value = 4110.5712
$ws.Range("M132").Value = -1580.5712
$ws.Range("H139").Value = 50715
$ws.Range("J139").Value = 50715
$ws.Range("L139").Value = 50715
$ws.Range("N139").Value = -60995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H134").Value = 2946.0889
$ws.Range("I134").Value = 3590.7742
$ws.Range("J134").Value = 1518.5714
$ws.Range("K134").Value = 10772.3226
$ws.Range("L134").Value = 4555.7142
$ws.Range("M134").Value = -8237.3226
$ws.Range("N134").Value = -9625.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3541.457
$ws.Range("I31").Value = 1979.4166
$ws.Range("J31").Value = 6949.5454
$ws.Range("K31").Value = 1979.4166
$ws.Range("L31").Value = 6949.5454
$ws.Range("M31").Value = -1684.4166
$ws.Range("N31").Value = -7539.5454
$ws.Range("H34").Value = 3541.457
$ws.Range("I34").Value = 1979.4166
$ws.Range("J34").Value = 6949.5454
$ws.Range("K34").Value = 1979.4166
$ws.Range("L34").Value = 6949.5454
$ws.Range("M34").Value = -1777.4166
$ws.Range("N34").Value = -7353.5454
$ws.Range("H122").Value = 5333.6665
$ws.Range("I122").Value = 5333.6665
$ws.Range("K122").Value = 16000.9995
$ws.Range("M122").Value = -13550.9995
$ws.Range("H132").Value = 2307.7144
$ws.Range("I132").Value = 1757.5294
$ws.Range("J132").Value = 21014
$ws.Range("K132").Value = 5272.5882
$ws.Range("L132").Value = 63042
$ws.Range("M132").Value = -2742.5882
$ws.Range("N132").Value = -68102
$ws.Range("H134").Value = 1120
$ws.Range("I134").Value = 1061.3334
$ws.Range("K134").Value = 3184.0002
$ws.Range("M134").Value = -649.0001999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1557.6666
$ws.Range("I5").Value = 1236.75
$ws.Range("J5").Value = 1814.4
$ws.Range("K5").Value = 3710.25
$ws.Range("L5").Value = 5443.200000000001
$ws.Range("M5").Value = -3598.25
$ws.Range("N5").Value = -5667.200000000001
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("M46").Value = -2909
$ws.Range("H50").Value = 148
$ws.Range("I50").Value = 122
$ws.Range("K50").Value = 366
$ws.Range("M50").Value = 115
$ws.Range("H53").Value = 148
$ws.Range("I53").Value = 122
$ws.Range("K53").Value = 366
$ws.Range("M53").Value = 115
$ws.Range("H86").Value = 55560784
$ws.Range("I86").Value = 678
$ws.Range("J86").Value = 125010920
$ws.Range("K86").Value = 2034
$ws.Range("L86").Value = 375032760
$ws.Range("M86").Value = -848
$ws.Range("N86").Value = -375035132
$ws.Range("H89").Value = 55560784
$ws.Range("I89").Value = 678
$ws.Range("J89").Value = 125010920
$ws.Range("K89").Value = 6102
$ws.Range("L89").Value = 1125098280
$ws.Range("M89").Value = -174
$ws.Range("N89").Value = -1125110136
$ws.Range("H131").Value = 233468.62
$ws.Range("I131").Value = 680
$ws.Range("J131").Value = 278733.1
$ws.Range("K131").Value = 2040
$ws.Range("L131").Value = 836199.2999999999
$ws.Range("M131").Value = 3000
$ws.Range("N131").Value = -846279.2999999999
$ws.Range("H135").Value = 1557.6666
$ws.Range("I135").Value = 1236.75
$ws.Range("J135").Value = 1814.4
$ws.Range("K135").Value = 11130.75
$ws.Range("L135").Value = 16329.6
$ws.Range("M135").Value = -8595.75
$ws.Range("N135").Value = -21399.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3537.6206
$ws.Range("I80").Value = 2700.5
$ws.Range("K80").Value = 2700.5
$ws.Range("M80").Value = -1702.5
$ws.Range("H83").Value = 3537.6206
$ws.Range("I83").Value = 2700.5
$ws.Range("K83").Value = 13502.5
$ws.Range("M83").Value = -8510.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1939.3572
$ws.Range("I22").Value = 1322.3334
$ws.Range("K22").Value = 1322.3334
$ws.Range("M22").Value = -1027.3334
$ws.Range("H27").Value = 1939.3572
$ws.Range("I27").Value = 1322.3334
$ws.Range("K27").Value = 1322.3334
$ws.Range("M27").Value = -1215.3334
$ws.Range("H46").Value = 1428.7
$ws.Range("J46").Value = 2075.5
$ws.Range("L46").Value = 2075.5
$ws.Range("N46").Value = -2451.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 891.8823
$ws.Range("I132").Value = 744.13336
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2232.40008
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = 297.5999199999997
$ws.Range("N132").Value = -11060
$ws.Range("H137").Value = 48043.332
$ws.Range("J137").Value = 48043.332
$ws.Range("L137").Value = 48043.332
$ws.Range("N137").Value = -58243.332
